$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2881169905109251
$ws.Range("C2").Value = 0.3048912486333797
$ws.Range("D2").Value = 0.1496068669990043
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("G2").Value = 1.276001064745008

$ws.Range("B3").Value = 0.01253208636536152
$ws.Range("C3").Value = 0.3048912486333797
$ws.Range("D3").Value = 18.71679738969934
$ws.Range("E3").Value = 13.86384647080068
$ws.Range("G3").Value = 32.89806719549876
